$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ozurgeti")

$ws.Range("E4").Value = 14218
$ws.Range("F4").Value = 14418
$ws.Range("G4").Value = 14513
$ws.Range("H4").Value = 14643
$ws.Range("I4").Value = 14751
$ws.Range("J4").Value = 14991
$ws.Range("K4").Value = 15066
